$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Methods Required")
$ws.Activate()

# Mark "Raw Code" as completed (x) for the Movement section rows (45-48)
$ws.Range("C45").Value = "x"

$ws.Range("C46").HorizontalAlignment = -4108
$ws.Range("C46").VerticalAlignment = -4108
$ws.Range("C46").Value = "x"

$ws.Range("C47").HorizontalAlignment = -4108
$ws.Range("C47").VerticalAlignment = -4108
$ws.Range("C47").Value = "x"

$ws.Range("C48").HorizontalAlignment = -4108
$ws.Range("C48").VerticalAlignment = -4108
$ws.Range("C48").Value = "x"

# New blank formatted row below the table body
$ws.Range("C49").HorizontalAlignment = -4108
$ws.Range("C49").VerticalAlignment = -4108

# Reflect final cursor/selection position as left by the author
$ws.Range("D50").Select()

Write-Host "Movement Raw Code rows updated"
